# 06_LibFormula.xlsx - "Indicatori BIB complessita 2 e 3 Matteo"
#
# - Re-sorts the existing INDICATOR_* rows (6-22) on sheet "Library_Formula"
#   into ascending numeric order.
# - Appends 20 new INDICATOR_* rows (23-42) for the newly added indicators.
# - Updates the view (active selection) and fixes up the "Normale" cell
#   style name to "Normal".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

# ---------------------------------------------------------------------------
# 1) Re-order the C column (Formula Name) for the existing indicator rows.
#    Only the <v> (shared-string) content changes for these rows - the
#    style ("s" attribute) of every cell here is already correct and must
#    stay untouched, so we only touch .Value here.
# ---------------------------------------------------------------------------
$existingValues = @{
    6  = "INDICATOR_14"
    7  = "INDICATOR_16"
    8  = "INDICATOR_20"
    9  = "INDICATOR_44"
    10 = "INDICATOR_49"
    11 = "INDICATOR_50"
    12 = "INDICATOR_55"
    13 = "INDICATOR_58"
    14 = "INDICATOR_193"
    15 = "INDICATOR_197"
    16 = "INDICATOR_210"
    17 = "INDICATOR_211 "
    18 = "INDICATOR_212"
    19 = "INDICATOR_213"
    20 = "INDICATOR_216"
    21 = "INDICATOR_217"
    22 = "INDICATOR_218"
}
foreach ($row in 6..22) {
    $ws.Cells.Item($row, 3).Value = $existingValues[$row]
}

# ---------------------------------------------------------------------------
# 2) Append the 20 new indicator rows (23-42), each following the same
#    A/B/C/E layout as the rows above it:
#      A = CREATE/MODIFY, B = LIB_EWS_BE, C = INDICATOR_xx, E = String
# ---------------------------------------------------------------------------
$newValues = @{
    23 = "INDICATOR_35"
    24 = "INDICATOR_56"
    25 = "INDICATOR_60"
    26 = "INDICATOR_66"
    27 = "INDICATOR_70"
    28 = "INDICATOR_76"
    29 = "INDICATOR_80"
    30 = "INDICATOR_84"
    31 = "INDICATOR_88"
    32 = "INDICATOR_92"
    33 = "INDICATOR_96"
    34 = "INDICATOR_100"
    35 = "INDICATOR_105"
    36 = "INDICATOR_113 "
    37 = "INDICATOR_118"
    38 = "INDICATOR_124"
    39 = "INDICATOR_128"
    40 = "INDICATOR_132"
    41 = "INDICATOR_173"
    42 = "INDICATOR_186"
}

# Match the existing font formatting (Trebuchet MS, 10pt, black) used by the
# rest of the table for columns A, B, C and E of the new rows.
$ws.Range("A23:A42").Font.Name = "Trebuchet MS"
$ws.Range("A23:A42").Font.Size = 10
$ws.Range("A23:A42").Font.Color = 0

$ws.Range("B23:B42").Font.Name = "Trebuchet MS"
$ws.Range("B23:B42").Font.Size = 10
$ws.Range("B23:B42").Font.Color = 0

$ws.Range("E23:E42").Font.Name = "Trebuchet MS"
$ws.Range("E23:E42").Font.Size = 10
$ws.Range("E23:E42").Font.Color = 0

$ws.Range("A23:A42,C23:C42").Font.Name = "Trebuchet MS"
$ws.Range("A23:A42,C23:C42").Font.Size = 10
$ws.Range("A23:A42,C23:C42").Font.Color = 0

foreach ($row in 23..42) {
    $ws.Cells.Item($row, 1).Value = "CREATE/MODIFY"
    $ws.Cells.Item($row, 2).Value = "LIB_EWS_BE"
    $ws.Cells.Item($row, 3).Value = $newValues[$row]
    $ws.Cells.Item($row, 5).Value = "String"
}

# ---------------------------------------------------------------------------
# 3) Update the view: the sheet now shows the bottom rows with C37 selected.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("C37").Select()

# ---------------------------------------------------------------------------
# 4) Rename the default cell style from "Normale" (Italian) to "Normal".
#    A plain rename is ignored for the builtin style, so remove/recreate it.
# ---------------------------------------------------------------------------
$wb.Styles.Item("Normale").Delete()
$wb.Styles.Add("Normal")
